# Data retrieved - Tue Jul 27 18:25:44 UTC 2021
# New daily snapshot appended as row 91; row 90's timestamp also re-saved
# with an (equal-value) floating point literal.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 90: refresh the date/time value (same instant, re-serialized literal).
$ws.Range("A90").Value = 44403.7681113912

# Row 91: new day's numbers appended to the bottom of the table.
$ws.Range("A91").Value = 44404.76787926738
$ws.Range("A91").NumberFormat = $ws.Range("A90").NumberFormat

$ws.Range("B91").Value = 80359
$ws.Range("C91").Value = 67809
$ws.Range("D91").Value = 3669
$ws.Range("E91").Value = 2216
$ws.Range("F91").Value = 1607
$ws.Range("G91").Value = 21063
$ws.Range("H91").Value = 1572
$ws.Range("I91").Value = 911
$ws.Range("J91").Value = 197
